$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (edit substrings of rich-text shared strings in place) ---
$ws.Range("A8").Characters(21,2).Text = "48"
$ws.Range("C9").Characters(27,10).Text = "11/28/2022"
$ws.Range("C9").Characters(48,10).Text = "12/4/2022"

# --- Crime data table updates (rows 15-30) ---
$ws.Range("F15").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = "'***.*"
$ws.Range("E15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("L15").Value = 16.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -28.571428571428
$ws.Range("I16").Value = 87
$ws.Range("J16").Value = 86
$ws.Range("K16").Value = 1.162790697674
$ws.Range("L16").Value = 52.631578947368
$ws.Range("M16").Value = 55.357142857142
$ws.Range("N16").Value = -85.051546391752
$ws.Range("C17").Value = 1
$ws.Range("D17").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -38.461538461538
$ws.Range("I17").Value = 106
$ws.Range("J17").Value = 111
$ws.Range("K17").Value = -4.504504504504
$ws.Range("L17").Value = 68.253968253968
$ws.Range("M17").Value = 107.843137254902
$ws.Range("N17").Value = -30.263157894736
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 170
$ws.Range("J18").Value = 141
$ws.Range("K18").Value = 20.567375886524
$ws.Range("L18").Value = -2.857142857142
$ws.Range("M18").Value = 95.402298850574
$ws.Range("N18").Value = -85.074626865671
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 128.571428571429
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 11.363636363636
$ws.Range("I19").Value = 612
$ws.Range("J19").Value = 499
$ws.Range("K19").Value = 22.645290581162
$ws.Range("L19").Value = 14.179104477611
$ws.Range("M19").Value = -10.263929618768
$ws.Range("N19").Value = -71.132075471698
$ws.Range("C20").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 200
$ws.Range("L20").Value = 32.608695652173
$ws.Range("N20").Value = -89.765100671140
$ws.Range("C21").Value = 21
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = 1.234567901234
$ws.Range("I21").Value = 1050
$ws.Range("J21").Value = 889
$ws.Range("K21").Value = 18.110236220472
$ws.Range("L21").Value = 18.110236220472
$ws.Range("M21").Value = 15.005476451259
$ws.Range("N21").Value = -77.183833116036
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 66.666666666666
$ws.Range("I22").Value = 30
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 42.857142857142
$ws.Range("M22").Value = -6.25
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 11.111111111111
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = -3.896103896103
$ws.Range("I24").Value = 1165
$ws.Range("J24").Value = 928
$ws.Range("K24").Value = 25.538793103448
$ws.Range("L24").Value = -4.742436631234
$ws.Range("M24").Value = 98.129251700680
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -6.25
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 167
$ws.Range("K25").Value = 31.736526946107
$ws.Range("L25").Value = 49.659863945578
$ws.Range("M25").Value = 2.803738317757
$ws.Range("F26").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = -100
$ws.Range("L26").Value = -11.111111111111
$ws.Range("C27").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = 23.809523809523
$ws.Range("F30").Value = "'0"
$ws.Range("E15").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$excel.CutCopyMode = 0
